$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (AVEIRO / MEALHADA) values for columns H through AA
$ws.Range("H2").Value = 208
$ws.Range("I2").Value = 524
$ws.Range("J2").Value = 2118
$ws.Range("L2").Value = 612
$ws.Range("M2").Value = 36
$ws.Range("N2").Value = 407
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 24
$ws.Range("S2").Value = 241
$ws.Range("T2").Value = 380
$ws.Range("U2").Value = 33
$ws.Range("V2").Value = 3479
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 3342
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 61
$ws.Range("AA2").Value = 16
